$d = $word.ActiveDocument

# Locate the unique occurrence of "MIEM Management Office" (typo for "NIEM")
$found = $d.Content.Find.Execute("MIEM Management Office", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)

$r = $d.Content
$r.Find.Execute("MIEM Management Office", $true, $false, $false, $false, $false,
                 $true, 1, $false, "", 0)

# Narrow the matched range down to just the first character ("M") and replace it with "N"
$r.SetRange($r.Start, $r.Start + 1)
$r.Text = "N"
